$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = 44322
Write-Host ("G9 value=" + $ws.Range("G9").Value)
